# Auto-generated cell updates reflecting the crypto price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.093.63'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '3.473.70'
$ws.Range('E3').Value = '  -0.46%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.26'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '177.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.22%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.594'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.82%  '
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').Value = '3.475.56'
$ws.Range('E9').Value = '  -0.39%  '
$ws.Range('E10').Value = '  +4.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.06'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.432'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.17%  '
$ws.Range('D13').Value = '4.076.37'
$ws.Range('E13').Value = '  -0.44%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '31.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +9.60%  '
$ws.Range('E15').Value = '  +0.89%  '
$ws.Range('D16').Value = '67.188.23'
$ws.Range('E16').Value = '  +0.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000176'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.89%  '
$ws.Range('D18').Value = '3.478.64'
$ws.Range('E18').Value = '  +0.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.61%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '386.84'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '73.86'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('E24').Value = '  -0.31%  '
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.534'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.31'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.173'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.91%  '
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.10'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.05'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.42'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.33'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.13%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.00%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.58'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.866'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.24%  '
$ws.Range('E40').Value = '  -1.79%  '
$ws.Range('E41').Value = '  +6.78%  '
$ws.Range('E42').Value = '  -2.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.62'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').Value = '2.826.56'
$ws.Range('E44').Value = '  +0.12%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '27.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0719'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.92%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.54'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.90%  '
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '333.60'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.76%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.62%  '
